$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Fix the wording in the "Built Environment..." paragraph:
#    "has facilitated a step towards" -> "have facilitated steps towards"
# -----------------------------------------------------------------
[void]$d.Content.Find.Execute(
    "has facilitated a step towards understanding",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "have facilitated steps towards understanding", 2)

# -----------------------------------------------------------------
# 2) Turn each of the five subsection titles (currently styled as
#    "First Paragraph") into their own "Heading 2" paragraph, and
#    turn the paragraph that follows them (currently styled
#    "Body Text") into a "First Paragraph" style, matching the other
#    subsections' pattern elsewhere in the document.
# -----------------------------------------------------------------
$titles = @(
    "Built Environment and Pediatric Psychiatric Disorders",
    "Causal Mediation of Place-Based Factors on Pediatric Health Disparities",
    "Fairness in Pediatric Precision Medicine",
    "Privacy-based Methods and Software for Geocoding and Geomarker Assessment",
    "Spatiotemporal Exposure Assessment Methods and Machine Learning Models"
)

$bodies = @(
    "Building on advanced exposure assessment",
    "I have applied advanced causal modeling techniques",
    "My research group has lead several studies",
    "Our group developed and maintains a novel approach",
    "My early career was spent developing"
)

foreach ($p in $d.Paragraphs) {
    $styleName = $p.Style.NameLocal
    $t = $p.Range.Text

    if ($styleName -eq "First Paragraph") {
        foreach ($title in $titles) {
            if ($t -like "$title*") {
                $p.Style = "Heading 2"
            }
        }
    } elseif ($styleName -eq "Body Text") {
        foreach ($body in $bodies) {
            if ($t -like "$body*") {
                $p.Style = "First Paragraph"
            }
        }
    }
}

# -----------------------------------------------------------------
# 3) Add bookmarks that wrap each subsection (heading + body +
#    citation list), splitting the previous single
#    "c.-contributions-to-science" bookmark into one bookmark per
#    subsection, while keeping the outer bookmark spanning the
#    whole "C. Contributions to Science" section.
# -----------------------------------------------------------------
$bookmarkNames = @(
    "X887e03c63b1cff1703d02ee410b927b2b453fee",
    "X85986e86d026ec2acdf2f0807db4d0bf595502a",
    "fairness-in-pediatric-precision-medicine",
    "Xd9163a628e6fdb35511b2a8702dc901b13c46e9",
    "X1918f807ef3509433f5daebb9f591f51071797f"
)

# Collect the paragraph index (1-based) of each Heading 2 title so we
# can compute the start/end of each subsection range afterwards.
$headingStarts = New-Object System.Collections.ArrayList
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $styleName = $p.Style.NameLocal
    $t = $p.Range.Text
    if ($styleName -eq "Heading 2") {
        foreach ($title in $titles) {
            if ($t -like "$title*") {
                [void]$headingStarts.Add($p.Range.Start)
            }
        }
    }
}

$docEnd = $d.Content.End

for ($k = 0; $k -lt $headingStarts.Count; $k++) {
    $start = $headingStarts[$k]
    if ($k -lt ($headingStarts.Count - 1)) {
        $end = $headingStarts[$k + 1]
    } else {
        $end = $docEnd
    }
    $r = $d.Range($start, $end)
    [void]$r.Bookmarks.Add($bookmarkNames[$k])
}

# Note: the outer "c.-contributions-to-science" bookmark already
# exists in the source document and already spans from right before
# the "C. Contributions to Science" heading to the end of the
# document/section, which is exactly the range it should keep. Only
# its numeric id changes on save (which the runtime renumbers
# automatically), so it needs no further edits here.
